$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Legislature" entity row (row 12) is being removed from the interaction matrix.
# Select the entire row and delete it, shifting all rows below up by one.
$ws.Rows.Item(12).Delete()
